$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.481.50'
$ws.Range("E2").Value = '  +0.32%  '

$ws.Range("D3").Value = '1.811.63'
$ws.Range("E3").Value = '  +1.26%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '228.60'
$ws.Range("E5").Value = '  +0.70%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.581'
$ws.Range("E6").Value = '  +4.56%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '35.90'
$ws.Range("E8").Value = '  +9.32%  '

$ws.Range("E9").Value = '  +2.83%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0699'
$ws.Range("E10").Value = '  +1.42%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0965'
$ws.Range("E11").Value = '  +1.98%  '

$ws.Range("D12").Value = '2.074.21'
$ws.Range("E12").Value = '  +1.17%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.56'
$ws.Range("E13").Value = '  +3.03%  '

$ws.Range("D14").Value = '1.790.81'
$ws.Range("E14").Value = '  +0.42%  '

$ws.Range("E15").Value = '  +2.02%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.55'
$ws.Range("E16").Value = '  +5.97%  '

$ws.Range("D17").Value = '34.478.31'
$ws.Range("E17").Value = '  +0.32%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.39'
$ws.Range("E18").Value = '  +1.11%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '247.83'
$ws.Range("E19").Value = '  +1.04%  '

$ws.Range("E20").Value = '  +0.46%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.57'
$ws.Range("E21").Value = '  +2.58%  '

$ws.Range("E22").Value = '  +0.09%  '

$ws.Range("E23").Value = '  +1.75%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '172.27'
$ws.Range("E24").Value = '  +2.30%  '

$ws.Range("E25").Value = '  +3.13%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.02'
$ws.Range("E26").Value = '  +9.06%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.04'
$ws.Range("E27").Value = '  +2.93%  '

$ws.Range("E28").Value = '  +3.80%  '

$ws.Range("E29").Value = '  -0.12%  '

$ws.Range("E30").Value = '  +1.93%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.90'
$ws.Range("E31").Value = '  +2.76%  '

$ws.Range("E32").Value = '  +1.83%  '

$ws.Range("E33").Value = '  +0.79%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.85'
$ws.Range("E34").Value = '  +1.84%  '

$ws.Range("D35").Value = '1.400.98'
$ws.Range("E35").Value = '  -0.88%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.681'
$ws.Range("E36").Value = '  -0.66%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.53'
$ws.Range("E37").Value = '  -1.63%  '

$ws.Range("E38").Value = '  +0.40%  '

$ws.Range("E39").Value = '  +0.79%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '83.79'
$ws.Range("E40").Value = '  -0.86%  '

$ws.Range("B41").Value = 'WEMIXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.24'
$ws.Range("E41").Value = '  +11.99%  '

$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.970'
$ws.Range("E42").Value = '  +2.84%  '

$ws.Range("E43").Value = '  +2.14%  '

$ws.Range("E44").Value = '  +0.25%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.44'
$ws.Range("E45").Value = '  -4.25%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.06'
$ws.Range("E46").Value = '  -1.33%  '

$ws.Range("E47").Value = '  -3.93%  '

$ws.Range("D48").Value = '1.973.63'
$ws.Range("E48").Value = '  +1.04%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '105.95'
$ws.Range("E49").Value = '  +0.67%  '

$ws.Range("E50").Value = '  +0.00%  '

$ws.Range("E51").Value = '  +0.98%  '
